# Adds rows 252-255 (new daily data through 2021-05-13) to Sheet1,
# mirroring the formatting of the preceding data row (A251).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 4,50
# row 252
$data[0,0] = 44326
$data[0,1] = 4
$data[0,2] = 0
$data[0,3] = 0
$data[0,4] = 10
$data[0,5] = 2
$data[0,6] = 1
$data[0,7] = 0
$data[0,8] = 0
$data[0,9] = 0
$data[0,10] = 0
$data[0,11] = 1
$data[0,12] = 0
$data[0,13] = 2
$data[0,14] = 0
$data[0,15] = 1
$data[0,16] = 0
$data[0,17] = 1
$data[0,18] = 0
$data[0,19] = 5
$data[0,20] = 0
$data[0,21] = 39
$data[0,22] = 0
$data[0,23] = 5
$data[0,24] = 3
$data[0,25] = 1
$data[0,26] = 0
$data[0,27] = 0
$data[0,28] = 1
$data[0,29] = 0
$data[0,30] = 0
$data[0,31] = 0
$data[0,32] = 9
$data[0,33] = 1
$data[0,34] = 0
$data[0,35] = 0
$data[0,36] = 1
$data[0,37] = 0
$data[0,38] = 3
$data[0,39] = 1
$data[0,40] = 2
$data[0,41] = 94
$data[0,42] = 0
$data[0,43] = 0
$data[0,44] = 0
$data[0,45] = 0
$data[0,46] = 0
$data[0,47] = 0
$data[0,48] = 1
$data[0,49] = 0

# row 253
$data[1,0] = 44327
$data[1,1] = 1
$data[1,2] = 0
$data[1,3] = 0
$data[1,4] = 10
$data[1,5] = 1
$data[1,6] = 0
$data[1,7] = 2
$data[1,8] = 0
$data[1,9] = 1
$data[1,10] = 1
$data[1,11] = 1
$data[1,12] = 1
$data[1,13] = 1
$data[1,14] = 0
$data[1,15] = 0
$data[1,16] = 0
$data[1,17] = 5
$data[1,18] = 1
$data[1,19] = 1
$data[1,20] = 1
$data[1,21] = 12
$data[1,22] = 0
$data[1,23] = 4
$data[1,24] = 1
$data[1,25] = 11
$data[1,26] = 0
$data[1,27] = 0
$data[1,28] = 0
$data[1,29] = 0
$data[1,30] = 0
$data[1,31] = 0
$data[1,32] = 4
$data[1,33] = 1
$data[1,34] = 2
$data[1,35] = 0
$data[1,36] = 1
$data[1,37] = 1
$data[1,38] = 2
$data[1,39] = 0
$data[1,40] = 3
$data[1,41] = 69
$data[1,42] = 0
$data[1,43] = 0
$data[1,44] = 0
$data[1,45] = 0
$data[1,46] = 0
$data[1,47] = 0
$data[1,48] = 0
$data[1,49] = 0

# row 254
$data[2,0] = 44328
$data[2,1] = 0
$data[2,2] = 0
$data[2,3] = 0
$data[2,4] = 0
$data[2,5] = 1
$data[2,6] = 0
$data[2,7] = 0
$data[2,8] = 0
$data[2,9] = 0
$data[2,10] = 0
$data[2,11] = 0
$data[2,12] = 0
$data[2,13] = 4
$data[2,14] = 0
$data[2,15] = 0
$data[2,16] = 1
$data[2,17] = 0
$data[2,18] = 0
$data[2,19] = 1
$data[2,20] = 1
$data[2,21] = 13
$data[2,22] = 0
$data[2,23] = 0
$data[2,24] = 1
$data[2,25] = 2
$data[2,26] = 0
$data[2,27] = 1
$data[2,28] = 1
$data[2,29] = 0
$data[2,30] = 0
$data[2,31] = 0
$data[2,32] = 0
$data[2,33] = 0
$data[2,34] = 0
$data[2,35] = 0
$data[2,36] = 0
$data[2,37] = 1
$data[2,38] = 0
$data[2,39] = 0
$data[2,40] = 1
$data[2,41] = 28
$data[2,42] = 0
$data[2,43] = 0
$data[2,44] = 0
$data[2,45] = 0
$data[2,46] = 0
$data[2,47] = 0
$data[2,48] = 0
$data[2,49] = 0

# row 255
$data[3,0] = 44329
$data[3,1] = 0
$data[3,2] = 0
$data[3,3] = 0
$data[3,4] = 14
$data[3,5] = 1
$data[3,6] = 4
$data[3,7] = 3
$data[3,8] = 0
$data[3,9] = 1
$data[3,10] = 0
$data[3,11] = 3
$data[3,12] = 1
$data[3,13] = 5
$data[3,14] = 0
$data[3,15] = 0
$data[3,16] = 0
$data[3,17] = 1
$data[3,18] = 0
$data[3,19] = 1
$data[3,20] = 2
$data[3,21] = 50
$data[3,22] = 1
$data[3,23] = 3
$data[3,24] = 1
$data[3,25] = 11
$data[3,26] = 0
$data[3,27] = 1
$data[3,28] = 0
$data[3,29] = 1
$data[3,30] = 1
$data[3,31] = 0
$data[3,32] = 6
$data[3,33] = 0
$data[3,34] = 0
$data[3,35] = 0
$data[3,36] = 1
$data[3,37] = 1
$data[3,38] = 2
$data[3,39] = 2
$data[3,40] = 1
$data[3,41] = 119
$data[3,42] = 0
$data[3,43] = 0
$data[3,44] = 0
$data[3,45] = 0
$data[3,46] = 0
$data[3,47] = 0
$data[3,48] = 1
$data[3,49] = 0

$ws.Range("A252:AX255").Value = $data

# Copy the date-column number format/border/alignment style (s="2") from
# the last existing row onto the newly added date cells in column A.
$ws.Range("A251").Copy() | Out-Null
$ws.Range("A252:A255").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
